# "Generate Report for Handback"
#
# For the 15a9ac79-6e20-493f-8ed3-72c13c6a1391 row (row 6) on both the
# zh-cn and de-de sheets, the handback-status checker found a handback
# file whose version doesn't match the latest commit, so it now:
#   - records the handback target (Latest Target File -> col I) as a
#     hyperlink to the handback .md, same as col A
#   - copies the handoff xlf name into the handback file column (col J)
#   - stamps the handback datetime (col K)
#   - writes the version-mismatch error detail (col P)
# Also the report widens columns I and P to fit the new content (40
# chars, matching the other wide columns on the sheet).

$wb = $excel.ActiveWorkbook

$handbackErrorMessage = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/067b3b21544316a3760dfc096c9b965ecf2fbffd/e2e/15a9ac79-6e20-493f-8ed3-72c13c6a1391.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/afa5556f3562ddc4d02b3592c54da93a5c3e44dd/e2e/15a9ac79-6e20-493f-8ed3-72c13c6a1391.md."

$sheetsInfo = @(
    @{ Name = "zh-cn"; HandbackDateTime = "2016-08-14 16:58:15" },
    @{ Name = "de-de"; HandbackDateTime = "2016-08-14 16:58:25" }
)

# Same handback-markdown target that col A already links to for this row
# (https://.../blob/afa5556f.../e2e/15a9ac79-....md) - identical on both
# sheets since the source markdown isn't per-language.
$mdAddress = "https://github.com/OpenLocalizationTestOrg/oltest/blob/afa5556f3562ddc4d02b3592c54da93a5c3e44dd/e2e/15a9ac79-6e20-493f-8ed3-72c13c6a1391.md"
$mdDisplay = "15a9ac79-6e20-493f-8ed3-72c13c6a1391.md"

foreach ($info in $sheetsInfo) {
    $ws = $wb.Worksheets.Item($info.Name)

    # Widen columns I (Latest Target File) and P (Error Detail) to 40,
    # matching the other full-width columns on the sheet.
    $ws.Range("I1").EntireColumn.ColumnWidth = 40
    $ws.Range("P1").EntireColumn.ColumnWidth = 40

    # Row 6 is the 15a9ac79-6e20-493f-8ed3-72c13c6a1391 file.
    $handoffFile = $ws.Range("G6").Value()

    # Col I ("Latest Target File"): hyperlink to the handback .md, same
    # target/display text as col A's hyperlink for this row.
    $ws.Hyperlinks.Add($ws.Range("I6"), $mdAddress, "", "", $mdDisplay) | Out-Null

    # Col J ("Latest Handback File"): the handed-back xlf, same as the
    # latest handoff file for this row.
    $ws.Range("J6").Value = $handoffFile

    # Col K ("Latest Handback DateTime").
    $ws.Range("K6").Value = $info.HandbackDateTime

    # Col P ("Error Detail"): version mismatch message.
    $ws.Range("P6").Value = $handbackErrorMessage
}
